# Adds a new "ODI Bowling Extra" worksheet (mirroring the existing
# "ODI Batting Extra" sheet) containing MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL
# scraping output, and tidies up the blank placeholder cells that the old
# "ODI Batting Extra" scraper used to leave behind.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clean up "ODI Batting Extra": the old scraper wrote an empty inline
#    string into every unused cell of columns B:E. The refreshed scraper
#    (this commit) only writes a cell when it actually has data, so those
#    placeholder cells should simply not exist any more.
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "B2","C2","D2","E2",
    "B5","C5","D5","E5",
    "B12","C12","D12","E12",
    "B14","C14","D14","E14",
    "B16","C16","D16","E16",
    "C17","D17","E17",
    "B18","C18","D18","E18",
    "B19","C19","D19","E19",
    "E20",
    "B21","C21","D21","E21"
)
foreach ($ref in $emptyCells) {
    $battingExtra.Range($ref).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet after "ODI Batting Extra" (i.e.
#    as the last / 5th tab).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Header row, styled like the other "Extra" sheet's header (bold, centered,
# top-aligned, thin box border).
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL). All three
# columns are text in the source data (including numeric-looking values),
# so force a text number format before writing so Excel doesn't coerce
# "0" / "30.00%" into numbers/percentages.
$data = @(
    @("4413","0",""),
    @("4417","0",""),
    @("4463","0","30.00%"),
    @("4464","",""),
    @("4465","0",""),
    @("4470","0",""),
    @("4480","0","20.00%"),
    @("4482","0",""),
    @("4487","0",""),
    @("4488","0","10.00%"),
    @("4491","0",""),
    @("4594","",""),
    @("4597","0","20.00%"),
    @("4600","",""),
    @("4601","0","20.00%"),
    @("4671","",""),
    @("4674","0","10.00%"),
    @("4675","",""),
    @("4687","",""),
    @("4689","0","10.00%")
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $r[2]
    $row++
}
